$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# "Analyze and design the REST API  with Endpoints and service layer"
#   -> collapse the double space between "API" and "with"
Replace-Text "Analyze and design the REST API  with Endpoints and service layer" "Analyze and design the REST API with Endpoints and service layer"

# "Implement DAO implementation  for existing customer to browse  loan Details"
#   -> collapse double spaces
Replace-Text "Implement DAO implementation  for existing customer to browse  loan Details" "Implement DAO implementation for existing customer to browse loan Details"

# "TEST the REST API  implementation by Junit" -> collapse double space
Replace-Text "REST API  implementation by Junit" "REST API implementation by Junit"

# "Perform Backend API Testing  through  POST man" (appears several times,
# including the one that starts mid-run with "Backend") -> collapse double spaces
Replace-Text "Testing  through  POST man" "Testing through POST man"

# "Clerk adding new customer and  loan application" -> collapse double space
Replace-Text "Clerk adding new customer and  loan application" "Clerk adding new customer and loan application"

# "Find all  loan application" -> collapse double space
Replace-Text "Find all  loan application" "Find all loan application"

# "Search by  date of apply" -> collapse double space
Replace-Text "Search by  date of apply" "Search by date of apply"

# "Existing customers can apply for topup loans" -> "top-up"
Replace-Text "Existing customers can apply for topup loans" "Existing customers can apply for top-up loans"

# "Creating  Roles " -> "Creating Roles "
Replace-Text "ing  Roles " "ing Roles "

# "Create Login Page  with roles" -> collapse double space
Replace-Text "Create Login Page  with roles" "Create Login Page with roles"

# "Manager to  search all loan application and documents" -> collapse double space
Replace-Text "Manager to  search all loan application and documents" "Manager to search all loan application and documents"
